$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (D3) and Correspond Handback DateTime (G3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 03:43:53"
$wsZhCn.Range("G3").Value = "2016-01-18 03:44:55"

# de-de sheet: Correspond Handoff Datetime (D3) and Correspond Handback DateTime (G3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 03:44:08"
$wsDeDe.Range("G3").Value = "2016-01-18 03:45:21"
